{"js": "// Fix: \"DB2\" was listed twice in the \"Used skills and tools\" list\n// (\"... Vagrant, MySQL, DB2, Jenkins, DB2\"). Remove the stray, duplicate\n// \", DB2\" that trails \"Jenkins\" while leaving the legitimate\n// \"MySQL, DB2\" mention untouched.\n\nconst body = context.document.body;\n\n// There are exactly two occurrences of the literal \", DB2\" in the\n// document: the legitimate \"MySQL, DB2\" earlier in the list, and the\n// duplicated \"Jenkins, DB2\" right after it. Search returns them in\n// document order, so the second hit is the one to remove.\nconst results = body.search(\", DB2\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length < 2) {\n  throw new Error(`Expected at least 2 matches of \", DB2\", found ${results.items.length}`);\n}\n\n// Delete just the duplicated \", DB2\" span (the second match) in place,\n// rather than rewriting/replacing the surrounding \"Jenkins\" text, so\n// neighboring bookmarks (e.g. the one ending right after \"Jenkins\")\n// are left completely intact.\nconst duplicate = results.items[1];\nduplicate.delete();\n\nawait context.sync();\n", "ps1": "# Fix: \"DB2\" was listed twice in the \"Used skills and tools\" list\n# (\"... Vagrant, MySQL, DB2, Jenkins, DB2\"). Remove the stray, duplicate\n# \", DB2\" that trails \"Jenkins\" while leaving the legitimate\n# \"MySQL, DB2\" mention untouched.\n\n$d = $word.ActiveDocument\n\n# There are exactly two occurrences of the literal \", DB2\" in the\n# document: the legitimate \"MySQL, DB2\" earlier in the list, and the\n# duplicated \"Jenkins, DB2\" right after it. Walk all matches (in\n# document order) and keep snapshots of their ranges.\n$matches = @()\n$searchRange = $d.Content\n$find = $searchRange.Find\n$find.ClearFormatting()\n$find.Text = \", DB2\"\n$find.Forward = $true\n$find.Wrap = 0\n$find.MatchCase = $true\nwhile ($find.Execute()) {\n  $matches += $d.Range($searchRange.Start, $searchRange.End)\n  $searchRange.Collapse(0)\n}\n\nif ($matches.Count -lt 2) {\n  throw \"Expected at least 2 matches of ', DB2', found $($matches.Count)\"\n}\n\n# Delete just the duplicated \", DB2\" span (the second match) in place,\n# rather than rewriting/replacing the surrounding \"Jenkins\" text, so\n# neighboring bookmarks (e.g. the one ending right after \"Jenkins\")\n# are left completely intact.\n$matches[1].Delete()\n"}
